# Add new columns I ("I0") and J ("IF") to Sheet1, mirroring the header style
# used by the other header cells (B1:H1), and filling data rows 2-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
# Re-use the style already applied to the existing header cells (H1 uses style index 1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$iValues = @(9,7,8,9,7,8,8,9,8,8,9,8,7,10,7,8,8,9,8,7,7,5,8,5,9,7,7,5,6,8,9,5,9,6,6,9,6,7,7,6,6,6,5,7,9,4,8,6,8,7,8,7,5,7,5)
$jValues = @(9,8,8,9,7,8,8,9,8,8,9,8,7,10,7,9,8,9,8,7,7,5,9,6,9,8,7,5,6,8,9,6,9,6,6,9,6,7,7,6,6,6,5,7,9,4,8,6,8,7,8,7,5,7,5)

for ($r = 2; $r -le 56; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
